$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-22 Thursday" "2024-02-23 Friday"

Replace-Text "726×7=5082" "617×4=2468"
Replace-Text "651×6=3906" "657×3=1971"
Replace-Text "929×4=3716" "434×8=3472"
Replace-Text "295×4=1180" "399×7=2793"
Replace-Text "645×7=4515" "573×3=1719"

Replace-Text "785×8=6280" "487×2=974"
Replace-Text "980×9=8820" "218×5=1090"
Replace-Text "764×7=5348" "429×6=2574"
Replace-Text "733×2=1466" "585×2=1170"
Replace-Text "112×4=448" "132×9=1188"

Replace-Text "948×3=2844" "448×8=3584"
Replace-Text "736×9=6624" "369×2=738"
Replace-Text "627×8=5016" "972×2=1944"
Replace-Text "645×8=5160" "435×4=1740"
Replace-Text "431×6=2586" "764×8=6112"

Replace-Text "948×8=7584" "574×9=5166"
Replace-Text "390×4=1560" "276×6=1656"
Replace-Text "418×5=2090" "982×3=2946"
Replace-Text "539×5=2695" "312×2=624"
Replace-Text "476×9=4284" "604×7=4228"

Replace-Text "906×5=4530" "293×3=879"
Replace-Text "510×3=1530" "965×2=1930"
Replace-Text "314×2=628" "185×2=370"
Replace-Text "773×7=5411" "692×7=4844"
Replace-Text "405×6=2430" "311×6=1866"
